$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  -2.55%  '
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.534.58'
$ws.Range("D2").Style = "Normal"

$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.672.30'
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("D4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"

$ws.Range("E5").Value = '  -1.46%  '
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.74'
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = '  -2.31%  '

$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.006'
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = '  -1.25%  '
$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06490'
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = '  -2.21%  '
$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2586'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  -3.17%  '
$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.11'
$ws.Range("D10").Style = "Normal"

$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07679'
$ws.Range("D11").Style = "Normal"

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.678.78'
$ws.Range("D12").Style = "Normal"

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E13").Value = '  -4.61%  '
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.351'
$ws.Range("D13").Style = "Normal"

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("E14").Value = '  -1.75%  '
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.902.68'
$ws.Range("D14").Style = "Normal"

$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5602'
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = '  -1.05%  '
$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8059'
$ws.Range("D16").Style = "Normal"

$ws.Range("E17").Value = '  -3.66%  '
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.97'
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = '  -2.38%  '
$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.555.16'
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = '  -1.41%  '
$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '211.80'
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = '  -4.29%  '
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.453'
$ws.Range("D21").Style = "Normal"

$ws.Range("E22").Value = '  -2.27%  '
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.15'
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.924'
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.007'
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.44'
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.733'
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1171'
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = '  -2.87%  '
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.028'
$ws.Range("D28").Style = "Normal"

$ws.Range("E29").Value = '  -2.85%  '
$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.79'
$ws.Range("D29").Style = "Normal"

$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05232'
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = '  -1.83%  '
$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.265'
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = '  -3.16%  '
$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.374'
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = '  -5.30%  '
$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.227'
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = '  -2.90%  '
$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.588'
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = '  -3.21%  '
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.775'
$ws.Range("D35").Style = "Normal"

$ws.Range("E36").Value = '  -1.95%  '
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.374'
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = '  -1.67%  '
$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9285'
$ws.Range("D37").Style = "Normal"

$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5746'
$ws.Range("D38").Style = "Normal"

$ws.Range("E39").Value = '  +11.63%  '
$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.161.44'
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01602'
$ws.Range("D40").Style = "Normal"

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E41").Value = '  +4.00%  '
$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8688'
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("E42").Value = '  +0.37%  '
$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.006'
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = '  -3.54%  '
$ws.Range("D43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.652'
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.35'
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = '  -1.79%  '
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.811.57'
$ws.Range("D45").Style = "Normal"

$ws.Range("E46").Value = '  +1.48%  '
$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈116'
$ws.Range("D46").Style = "Normal"

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E47").Value = '  -2.90%  '
$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.17'
$ws.Range("D47").Style = "Normal"

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4492'
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = '  +0.48%  '

$ws.Range("E50").Value = '  -1.10%  '
$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.988'
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = '  -1.83%  '
$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05140'
$ws.Range("D51").Style = "Normal"
